# Update cell values in DOT_hemo1 dataset (Sheet1) per source edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 93.422325043864817
$ws.Range("B2").Value = 89.75599155294401
$ws.Range("C2").Value = 77.06761800207245
$ws.Range("D2").Value = 92.89159298520751
$ws.Range("E2").Value = 88.449563637334634
$ws.Range("F2").Value = 92.133922026468852
$ws.Range("G2").Value = 94.060429651605887
$ws.Range("H2").Value = 99.319244882157207
$ws.Range("I2").Value = 100.79116490368604
$ws.Range("J2").Value = 98.538979573807396
$ws.Range("K2").Value = 96.343327645167562
$ws.Range("L2").Value = 96.352824085283203
$ws.Range("M2").Value = 93.763700113525388
$ws.Range("N2").Value = 96.288294729206811
$ws.Range("O2").Value = 92.155896825640369

# Row 3
$ws.Range("A3").Value = 84.410955116312536
$ws.Range("B3").Value = 84.980044169608561
$ws.Range("C3").Value = 84.892962167637847
$ws.Range("D3").Value = 84.119296063060801
$ws.Range("E3").Value = 81.950685961203931
$ws.Range("F3").Value = 80.636148785626929
$ws.Range("G3").Value = 86.129849659662227
$ws.Range("H3").Value = 86.575842454432092
$ws.Range("I3").Value = 86.869070257475357
$ws.Range("J3").Value = 86.467839812850585
$ws.Range("K3").Value = 85.210990466618227
$ws.Range("L3").Value = 82.782085057671594
$ws.Range("M3").Value = 83.321613885378952
$ws.Range("N3").Value = 83.307252191624116
$ws.Range("O3").Value = 84.509690695524512

# Row 37
$ws.Range("A37").Value = 94.48097819510545
$ws.Range("B37").Value = 94.48097819510545
$ws.Range("C37:I37").ClearContents()

# Row 72
$ws.Range("A72").Value = 78.762445800138238
$ws.Range("B72").Value = 73.216551677283277
$ws.Range("C72").Value = 82.823721938382533
$ws.Range("D72").Value = 78.376648595400525
$ws.Range("E72").Value = 79.482256486819409
$ws.Range("F72").Value = 79.015170564324691
$ws.Range("G72").Value = 81.908619796671488
$ws.Range("H72").Value = 78.633430139666046
$ws.Range("I72").Value = 76.643167202557933

# Row 99
$ws.Range("A99").Value = 102.1879871806351
$ws.Range("B99").Value = 99.273252401867197
$ws.Range("C99").Value = 102.94670584891755
$ws.Range("D99").Value = 100.81354095372939
$ws.Range("E99").Value = 107.74235068693692
$ws.Range("F99").Value = 107.33293116046683
$ws.Range("G99").Value = 95.01914203189277

# Row 145
$ws.Range("A145").Value = 79.607376614078248
$ws.Range("B145").Value = 78.452992997065834
$ws.Range("C145").Value = 79.370936261222212
$ws.Range("D145").Value = 82.314730720007688
$ws.Range("E145").Value = 78.872305113372377
$ws.Range("F145").Value = 78.657189768593369
$ws.Range("G145").Value = 79.094799427178941
$ws.Range("H145").Value = 80.032482406943998
$ws.Range("I145").Value = 79.747366398041635
$ws.Range("J145").Value = 79.301625763411181
$ws.Range("K145").Value = 80.229337284945302

# Row 179
$ws.Range("A179").Value = 52.033332622334058
$ws.Range("B179").Value = 55.60332467579412
$ws.Range("C179").Value = 55.214124902571385
$ws.Range("D179").Value = 48.556026108449281
$ws.Range("E179").Value = 48.759854802521424
$ws.Range("F179:N179").ClearContents()

# Update the active selection to match the final cursor position (A311:E311),
# scrolled so row 289 is near the top of the view.
$ws.Range("A311:E311").Select() | Out-Null
